# Fruta / hortaliza, semanal
# Update the weekly Chirimoya price records: refresh existing rows 2-6 with
# newer reported figures and append three new rows (8-10) for additional
# market days.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---------------------------------------------------------------
$ws.Range("D2").Value = 44445
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 250
$ws.Range("N2").Value = 28000
$ws.Range("O2").Value = 30000
$ws.Range("P2").Value = 29200
$ws.Range("S2").Value = 2920

# --- Row 3 ---------------------------------------------------------------
$ws.Range("D3").Value = 44431
$ws.Range("N3").Value = 30000
$ws.Range("O3").Value = 30000
$ws.Range("P3").Value = 30000
$ws.Range("S3").Value = 3000

# --- Row 4 ---------------------------------------------------------------
$ws.Range("D4").Value = 44421
$ws.Range("M4").Value = 30
$ws.Range("N4").Value = 35000
$ws.Range("O4").Value = 35000
$ws.Range("P4").Value = 35000
$ws.Range("S4").Value = 3500

# --- Row 5 ---------------------------------------------------------------
$ws.Range("D5").Value = 44434
$ws.Range("M5").Value = 60

# --- Row 6 ---------------------------------------------------------------
$ws.Range("D6").Value = 44441
$ws.Range("M6").Value = 150

# --- Row 7 is unchanged ----------------------------------------------------

# --- New row 8 -------------------------------------------------------------
$ws.Range("A8").Value = 5
$ws.Range("B8").Value = "Macroferia Regional de Talca"
$ws.Range("C8").Value = "Maule"
$ws.Range("D8").Value = 44446
$ws.Range("D8").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E8").Value = 7
$ws.Range("F8").Value = "Fruta"
$ws.Range("G8").Value = 100107
$ws.Range("H8").Value = "Otros"
$ws.Range("I8").Value = 100107002
$ws.Range("J8").Value = "Chirimoya"
$ws.Range("K8").Value = "Cultivar IV Región"
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 200
$ws.Range("N8").Value = 30000
$ws.Range("O8").Value = 30000
$ws.Range("P8").Value = 30000
$ws.Range("Q8").Value = "$/bandeja 10 kilos"
$ws.Range("R8").Value = "Provincia de Limarí"
$ws.Range("S8").Value = 3000
$ws.Range("T8").Value = 10

# --- New row 9 -------------------------------------------------------------
$ws.Range("A9").Value = 5
$ws.Range("B9").Value = "Macroferia Regional de Talca"
$ws.Range("C9").Value = "Maule"
$ws.Range("D9").Value = 44438
$ws.Range("D9").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E9").Value = 7
$ws.Range("F9").Value = "Fruta"
$ws.Range("G9").Value = 100107
$ws.Range("H9").Value = "Otros"
$ws.Range("I9").Value = 100107002
$ws.Range("J9").Value = "Chirimoya"
$ws.Range("K9").Value = "Cultivar IV Región"
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 100
$ws.Range("N9").Value = 30000
$ws.Range("O9").Value = 30000
$ws.Range("P9").Value = 30000
$ws.Range("Q9").Value = "$/bandeja 10 kilos"
$ws.Range("R9").Value = "Provincia de Limarí"
$ws.Range("S9").Value = 3000
$ws.Range("T9").Value = 10

# --- New row 10 ------------------------------------------------------------
$ws.Range("A10").Value = 5
$ws.Range("B10").Value = "Macroferia Regional de Talca"
$ws.Range("C10").Value = "Maule"
$ws.Range("D10").Value = 44435
$ws.Range("D10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E10").Value = 7
$ws.Range("F10").Value = "Fruta"
$ws.Range("G10").Value = 100107
$ws.Range("H10").Value = "Otros"
$ws.Range("I10").Value = 100107002
$ws.Range("J10").Value = "Chirimoya"
$ws.Range("K10").Value = "Cultivar IV Región"
$ws.Range("L10").Value = "Especial"
$ws.Range("M10").Value = 160
$ws.Range("N10").Value = 30000
$ws.Range("O10").Value = 30000
$ws.Range("P10").Value = 30000
$ws.Range("Q10").Value = "$/bandeja 10 kilos"
$ws.Range("R10").Value = "Provincia de Limarí"
$ws.Range("S10").Value = 3000
$ws.Range("T10").Value = 10
